# Auto update: 2025-11-29 03:48:30
# Update "최종점수" (K) and "MACRO_SCORE" (N) values for rows 2-5

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("K2").Value = 60.7
$ws.Range("N2").Value = 85.82376350509293

$ws.Range("K3").Value = 54.9
$ws.Range("N3").Value = 85.82376350509293

$ws.Range("K4").Value = 51.1
$ws.Range("N4").Value = 85.82376350509293

$ws.Range("K5").Value = 48.7
$ws.Range("N5").Value = 85.82376350509293
